$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '60.909.19'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -0.05%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.913.07'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +0.06%  '

$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '592.37'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +1.10%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '145.59'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -0.61%  '

$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.01%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.506'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.67%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '6.89'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +1.98%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.143'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -0.62%  '

$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -1.99%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0000225'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +0.29%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '33.50'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -0.30%  '

$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -0.46%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.399.41'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +0.19%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '60.918.17'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +0.09%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '6.68'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -1.27%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.916.75'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +0.25%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '429.96'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +0.74%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.35'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -1.75%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.678'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +0.86%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '7.05'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -0.96%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '81.46'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +1.51%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '10.95'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -0.23%  '

$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -0.32%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.86'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -0.12%  '

$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +5.39%  '

$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +0.04%  '

$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -0.54%  '

$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -2.84%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '26.48'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +0.05%  '

$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +0.79%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0₃0851'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +0.69%  '

$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +0.07%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.62'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -0.19%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.00'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +0.71%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.122'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -0.89%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.98'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -1.66%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '8.53'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -1.61%  '

$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -1.65%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '39.99'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -3.69%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '374.45'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -0.70%  '

$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.700.83'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +1.13%  '

$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0344'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -0.94%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '132.38'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -0.38%  '

$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -0.10%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '23.76'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -3.81%  '

$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -0.36%  '

$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -3.57%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.124'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +1.20%  '
